$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = 'Volume 31   Number  44'
$ws.Range("C9").Value = 'Report Covering the Week  10/28/2024  Through  11/3/2024'

# --- Cells changing from numeric to text placeholder ("0" / "***.*") ---
# Donor cells with style 13 + matching text, stable across this edit:
#   "0"    donor -> C14   (or F14/D27/...)
#   "***.*" donor -> E27  (or E15/H15/...)
$ws.Range("C14").Copy($ws.Range("D14"))
$ws.Range("E27").Copy($ws.Range("E14"))
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("E27").Copy($ws.Range("E29"))
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("E27").Copy($ws.Range("E30"))
$ws.Range("C14").Copy($ws.Range("C31"))

# --- Cells changing from text placeholder to numeric ---
# Donor cells with style 14 (int) / 15 (pct) stable across this edit:
$ws.Range("G14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("H14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("G14").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 2
$ws.Range("H14").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("G14").Copy($ws.Range("C33"))
$ws.Range("C33").Value = 1
$ws.Range("G14").Copy($ws.Range("D33"))
$ws.Range("D33").Value = 1
$ws.Range("H14").Copy($ws.Range("E33"))
$ws.Range("E33").Value = 0
$ws.Range("G14").Copy($ws.Range("F33"))
$ws.Range("F33").Value = 1

# --- Plain numeric value updates ---
$ws.Range("M14").Value = -75
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 5
$ws.Range("I15").Value = 27
$ws.Range("K15").Value = 42.105263157894
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 58.823529411764
$ws.Range("N15").Value = 35
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -22.727272727272
$ws.Range("I16").Value = 177
$ws.Range("J16").Value = 198
$ws.Range("K16").Value = -10.60606060606
$ws.Range("L16").Value = 4.117647058823
$ws.Range("M16").Value = -20.982142857142
$ws.Range("N16").Value = -78.75150060024
$ws.Range("C17").Value = 8
$ws.Range("E17").Value = 0
$ws.Range("I17").Value = 253
$ws.Range("J17").Value = 238
$ws.Range("K17").Value = 6.302521008403
$ws.Range("L17").Value = -1.556420233463
$ws.Range("M17").Value = 37.5
$ws.Range("N17").Value = -2.692307692307
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -18.75
$ws.Range("I18").Value = 187
$ws.Range("J18").Value = 171
$ws.Range("K18").Value = 9.356725146198
$ws.Range("L18").Value = -19.742489270386
$ws.Range("M18").Value = -50.918635170603
$ws.Range("N18").Value = -89.006466784244
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 62
$ws.Range("H19").Value = -29.032258064516
$ws.Range("I19").Value = 548
$ws.Range("J19").Value = 591
$ws.Range("K19").Value = -7.275803722504
$ws.Range("L19").Value = -1.438848920863
$ws.Range("M19").Value = 50.136986301369
$ws.Range("N19").Value = 1.293900184842
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 29
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = 31.818181818181
$ws.Range("I20").Value = 378
$ws.Range("J20").Value = 320
$ws.Range("K20").Value = 18.125
$ws.Range("L20").Value = 47.081712062256
$ws.Range("M20").Value = 20.766773162939
$ws.Range("N20").Value = -87.3281930942
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 131
$ws.Range("G21").Value = 145
$ws.Range("H21").Value = -9.655172413793
$ws.Range("I21").Value = 1571
$ws.Range("J21").Value = 1541
$ws.Range("K21").Value = 1.946787800129
$ws.Range("L21").Value = 5.153949129852
$ws.Range("M21").Value = 5.577956989247
$ws.Range("N21").Value = -75.279307631786
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = -20
$ws.Range("M22").Value = -5.882352941176
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 84
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = -15.151515151515
$ws.Range("I24").Value = 1180
$ws.Range("J24").Value = 1117
$ws.Range("K24").Value = 5.640107430617
$ws.Range("L24").Value = -4.684975767366
$ws.Range("M24").Value = 32.286995515695
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = -27.450980392156
$ws.Range("I25").Value = 464
$ws.Range("J25").Value = 395
$ws.Range("K25").Value = 17.468354430379
$ws.Range("L25").Value = -7.014028056112
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 15
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = 22.5
$ws.Range("I26").Value = 521
$ws.Range("J26").Value = 402
$ws.Range("K26").Value = 29.601990049751
$ws.Range("L26").Value = 13.260869565217
$ws.Range("M26").Value = -14.590163934426
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 34
$ws.Range("K27").Value = 13.333333333333
$ws.Range("L27").Value = 17.241379310344
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("J28").Value = 54
$ws.Range("K28").Value = -37.037037037037
$ws.Range("L28").Value = -46.031746031746
$ws.Range("M29").Value = -50
$ws.Range("M30").Value = -33.333333333333
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = -66.666666666666
$ws.Range("J31").Value = 3
$ws.Range("K31").Value = 133.333333333333
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = -50
$ws.Range("I33").Value = 3
$ws.Range("J33").Value = 8
$ws.Range("K33").Value = -62.5
$ws.Range("L33").Value = -25
